# complexityAnalysis finished: add the 7th test-run's second/third sub-results
# (row 27 / row 28) and fix up the "Laufzeit Write" text for the first
# sub-result (row 26 / K26). Also switches the little "Fails" summary cell
# F35 (and its new neighbour G35) from the stray date-time format to the
# elapsed-time duration format used elsewhere in the sheet, and moves the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27: fill in the results for sub-test "b" --------------------------
$ws.Range("G27").Value = "02:31:11h"
$ws.Range("H27").Value = "04:05m"
$ws.Range("I27").Value = "00:51:01h"
$ws.Range("J27").Value = "02:26:05h"
$ws.Range("K27").Value = "01:48:12h"
$ws.Range("L27").Value = "10,9GB"
$ws.Range("M27").Value = 1

# --- Row 26: correct the "Laufzeit Write" value for sub-test "a" ----------
$ws.Range("K26").Value = "01:40:51h"

# --- Row 28: fill in the results for sub-test "c" --------------------------
$ws.Range("G28").Value = "02:31:07h"
$ws.Range("H28").Value = "03:35m"
$ws.Range("I28").Value = "00:40:34h"
$ws.Range("J28").Value = "02:27:22h"
$ws.Range("K28").Value = "02:39:05h"
$ws.Range("L28").Value = "10,9GB"
$ws.Range("M28").Value = 0

# --- Summary block: give F35/G35 the elapsed-time format -------------------
$ws.Range("F35").NumberFormat = "[h]:mm:ss"
$ws.Range("G35").NumberFormat = "[h]:mm:ss"

# --- Move the active selection from D35 to D34 ------------------------------
$ws.Range("D34").Select()
